# Regenerate the 15 lattice-multiplication practice cells (5 rows x 3 cols)
# with a new set of problems, per the target revision.
#
# Each cell's paragraph is a single run containing 5 text segments separated
# by line breaks:
#   "{A} x {B}"
#   "  {tens(B)}    {ones(B)}"
#   "  ----"
#   "{tens(A)}|    |"
#   "{ones(A)}|    |"
#
# We rebuild each cell's text (preserving the existing run/paragraph
# formatting, e.g. sz=32) using the new problem list, in row-major order.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newProblems = @(
    "53 x 78", "16 x 46", "78 x 68",
    "41 x 77", "16 x 11", "39 x 38",
    "39 x 13", "51 x 50", "22 x 90",
    "67 x 65", "15 x 95", "92 x 59",
    "89 x 68", "38 x 83", "13 x 80"
)

$vtab = [char]11
$rows = 5
$cols = 3
$i = 0

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $problem = $newProblems[$i]
        $parts = $problem.Split("x")
        $a = $parts[0].Trim()
        $b = $parts[1].Trim()
        $b0 = $b.Substring(0,1)
        $b1 = $b.Substring(1,1)
        $a0 = $a.Substring(0,1)
        $a1 = $a.Substring(1,1)

        # NOTE: use string interpolation (not the "+" operator) to join these
        # pieces -- when both operands of "+" look like pure integers this
        # runtime performs numeric addition instead of concatenation.
        $factorLine = "  $b0    $b1"
        $digit1 = "$a0|    |"
        $digit2 = "$a1|    |"

        $newText = "$problem$vtab$factorLine$vtab  ----$vtab$digit1$vtab$digit2"

        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        # Exclude the trailing cell-mark/paragraph-mark characters
        $rng.End = $rng.End - 2
        $rng.Text = $newText

        $i = $i + 1
    }
}

Write-Host "Updated $i cells"
